$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A72").Value = "What's the maximum number of lithology types in an log?"
$ws.Range("B72").Value = "llama3.2:latest"
$ws.Range("C72").Value = "The maximum number of lithology types that can be recorded in a log is 450."
$ws.Range("A73").Value = "What's the maximum number of lithology types in an log?"
$ws.Range("B73").Value = "llama3.2:latest"
$ws.Range("C73").Value = "The maximum number of lithology types that can be recorded in a log is 450."
$ws.Range("A74").Value = "How many symbols can I have in the plot at any one time?"
$ws.Range("B74").Value = "llama3.2:latest"
$ws.Range("C74").Value = "You can have up to 10,000 symbols in a plot at any given time."
$ws.Range("A75").Value = "How many tables can I have in my log?"
$ws.Range("B75").Value = "llama3.2:latest"
$ws.Range("C75").Value = "You can have up to 100 tables in a log."
$ws.Range("A76").Value = "How many symbols can I have in the plot at any one time?"
$ws.Range("B76").Value = "llama3.2:latest"
$ws.Range("C76").Value = "You can have up to 10,000 symbols in a plot at any given time."
$ws.Range("A77").Value = "How many curves can I load in one go?"
$ws.Range("B77").Value = "llama3.2:latest"
$ws.Range("C77").Value = "You can load up to 450 curves at a time."
$ws.Range("A78").Value = "How many curves can I load in one go?"
$ws.Range("B78").Value = "llama3.2:latest"
$ws.Range("C78").Value = "You can load up to 450 curves at a time."
$ws.Range("A79").Value = "Can I improve the format"
$ws.Range("B79").Value = "llama3.2:latest"
$ws.Range("C79").Value = "Yes, you can improve the format of your text entries by using the Format Text context tab or the Rich Edit control. `nFor track text and free format text entries, you can use either option. If you want to align the text, you can use either the Format Text context tab or the Rich Edit control.`nTo start, highlight the text you wish to edit and then select the desired format from the available options in the Format Text context tab or the Rich Edit control."

Write-Host "New dimension / used range:" $ws.UsedRange.Address()
